$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the Storypoints value for the sprint (C4) from 10 to 7
$ws.Range("C4").Value = 7

# Update selection/scroll position to match the new view state
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 4

$wb.Save()
